$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-round a double to 15 significant decimal digits (the precision Excel
# itself stores/displays numbers at), so arithmetic results match the way
# Excel would have recomputed + redisplayed the value.
function Round15([double]$x) {
    if ($x -eq 0) { return 0.0 }
    return [double]($x.ToString("G15"))
}

# The B:F columns (B5_O, B5_C, B5_E, B5_A, B5_N) hold means of up to three
# underlying item scores, so every value is an exact multiple of 1/6. The
# values already on the sheet are themselves displayed/rounded to 15
# significant digits, so to faithfully halve the *true* underlying mean
# (rather than the already-rounded display value) we first recover the
# exact sixths numerator, halve that exactly, and then re-round the result
# to 15 significant digits the same way Excel would.
function HalveMeanOfThree([double]$value) {
    $sixths = [Math]::Round($value * 6)
    $halved = $sixths / 12
    return Round15 $halved
}

for ($row = 2; $row -le 28; $row++) {
    for ($col = 2; $col -le 6; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $current = $cell.Value2
        $cell.Value2 = HalveMeanOfThree $current
    }
}
